$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to plain text so numeric-looking strings (e.g. "29.019.98",
# "0.9980", "  -0.06%  ") are preserved literally as text, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.019.98"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "1.832.06"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "242.47"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "0.6266"
$ws.Range("E6").Value = "  -4.16%  "

$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.07604"
$ws.Range("E8").Value = "  +3.26%  "

$ws.Range("D9").Value = "0.2922"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").Value = "22.53"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").Value = "0.07716"
$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("D12").Value = "1.835.47"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").Value = "4.955"
$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").Value = "0.6636"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").Value = "0.00001020"
$ws.Range("E15").Value = "  +18.25%  "

$ws.Range("D16").Value = "82.71"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "6.047"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "29.020.63"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").Value = "226.30"
$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "7.180"
$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "158.52"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").Value = "8.488"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").Value = "0.1376"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").Value = "17.88"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").Value = "1.491"
$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("D29").Value = "4.094"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").Value = "4.003"
$ws.Range("E30").Value = "  -0.29%  "

$ws.Range("D31").Value = "1.186"
$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").Value = "0.05235"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("D33").Value = "1.843"
$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("D34").Value = "0.7333"
$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -1.39%  "

$ws.Range("D36").Value = "2.698"
$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("D37").Value = "1.236.82"
$ws.Range("E37").Value = "  -4.24%  "

$ws.Range("D38").Value = "2.756"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("D39").Value = "0.01783"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").Value = "6.330"
$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("D41").Value = "0.8967"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "101.82"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.975.38"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000125"
$ws.Range("E45").Value = "  +3.96%  "

$ws.Range("D46").Value = "64.16"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "0.5102"
$ws.Range("E47").Value = "  -0.75%  "

$ws.Range("D48").Value = "0.4036"
$ws.Range("E48").Value = "  +1.14%  "

$ws.Range("D49").Value = "8.857"
$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("D50").Value = "0.05738"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").Value = "6.672"
$ws.Range("E51").Value = "  -0.55%  "

# Restore default (General) formatting/style on these cells so only the
# cell contents change, matching the original styling.
$ws.Range("D2:E51").ClearFormats()
